$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.9586238692716051
$ws.Range("D2").Value = 0.009144461530631531
$ws.Range("E2").Value = 0.6722825207559993
$ws.Range("F2").Value = 0.384057915175795
$ws.Range("G2").Value = 0.2493200085657747
$ws.Range("H2").Value = 0.3688069070791613
$ws.Range("L2").Value = 0.1544792843802014
$ws.Range("M2").Value = 0.2034375728087099
$ws.Range("N2").Value = 1.784511145164188
$ws.Range("O2").Value = 1.152264400375913
# Row 3
$ws.Range("B3").Value = 0.9232221726074954
$ws.Range("D3").Value = 0.007958255668356173
$ws.Range("E3").Value = 0.6684553001518125
$ws.Range("F3").Value = 0.3700419434593272
$ws.Range("G3").Value = 0.2360756146370449
$ws.Range("H3").Value = 0.3659973866940902
$ws.Range("L3").Value = 0.138264602842483
$ws.Range("M3").Value = 0.1912961617525397
$ws.Range("N3").Value = 1.772666014063461
$ws.Range("O3").Value = 1.118026951725682
# Row 4
$ws.Range("B4").Value = 0.9018086496192268
$ws.Range("D4").Value = 0.007228256460006577
$ws.Range("E4").Value = 0.6662555964438113
$ws.Range("F4").Value = 0.3617052620121228
$ws.Range("G4").Value = 0.2281328915768626
$ws.Range("H4").Value = 0.3644839444711607
$ws.Range("L4").Value = 0.1283049262639651
$ws.Range("M4").Value = 0.1838900102612087
$ws.Range("N4").Value = 1.766078413195842
$ws.Range("O4").Value = 1.097839114247421
# Row 5
$ws.Range("B5").Value = 0.8931645848557821
$ws.Range("D5").Value = 0.006930374724575472
$ws.Range("E5").Value = 0.6653974009233821
$ws.Range("F5").Value = 0.3583756640998104
$ws.Range("G5").Value = 0.2249436815688171
$ws.Range("H5").Value = 0.3639204508781688
$ws.Range("L5").Value = 0.1242456165055899
$ws.Range("M5").Value = 0.1808844404367065
$ws.Range("N5").Value = 1.763567264345355
$ws.Range("O5").Value = 1.089821979151111
# Row 6
$ws.Range("B6").Value = 0.8917342288413295
$ws.Range("D6").Value = 0.006880887970279304
$ws.Range("E6").Value = 0.6652572180304901
$ws.Range("F6").Value = 0.3578268732493655
$ws.Range("G6").Value = 0.2244169822080977
$ws.Range("H6").Value = 0.3638301001561075
$ws.Range("L6").Value = 0.1235715399681112
$ws.Range("M6").Value = 0.1803861298627929
$ws.Range("N6").Value = 1.763160789950348
$ws.Range("O6").Value = 1.088503392742496
# Row 7
$ws.Range("B7").Value = 0.9016917391358561
$ws.Range("D7").Value = 0.00722424072483463
$ws.Range("E7").Value = 0.6662438672750639
$ws.Range("F7").Value = 0.3616600839357318
$ws.Range("G7").Value = 0.2280896885382617
$ws.Range("H7").Value = 0.3644761293625152
$ws.Range("L7").Value = 0.128250183277828
$ws.Range("M7").Value = 0.1838494251705214
$ws.Range("N7").Value = 1.766043843718521
$ws.Range("O7").Value = 1.097730144033875
# Row 8
$ws.Range("B8").Value = 0.9463507579103521
$ws.Range("D8").Value = 0.008735813160249961
$ws.Range("E8").Value = 0.6709319691538838
$ws.Range("F8").Value = 0.3791693220809194
$ws.Range("G8").Value = 0.2447139641292466
$ws.Range("H8").Value = 0.3677942824806308
$ws.Range("L8").Value = 0.1488894214248688
$ws.Range("M8").Value = 0.199241244492832
$ws.Range("N8").Value = 1.780285295912677
$ws.Range("O8").Value = 1.140286098249589
# Row 9
$ws.Range("B9").Value = 1.036460340040406
$ws.Range("D9").Value = 0.01168617385914672
$ws.Range("E9").Value = 0.6812993445192035
$ws.Range("F9").Value = 0.4156439479334182
$ws.Range("G9").Value = 0.278823685091993
$ws.Range("H9").Value = 0.3759795815469857
$ws.Range("L9").Value = 0.189321817104755
$ws.Range("M9").Value = 0.2298019763489663
$ws.Range("N9").Value = 1.813609527625786
$ws.Range("O9").Value = 1.230369855922874
# Row 10
$ws.Range("B10").Value = 1.104173268278544
$ws.Range("D10").Value = 0.01384473253049379
$ws.Range("E10").Value = 0.6896089852287801
$ws.Range("F10").Value = 0.4437536655391199
$ws.Range("G10").Value = 0.3048166408657238
$ws.Range("H10").Value = 0.3830168937690956
$ws.Range("L10").Value = 0.2189905900715274
$ws.Range("M10").Value = 0.2524750201102606
$ws.Range("N10").Value = 1.841332301050727
$ws.Range("O10").Value = 1.300624483796497
# Row 11
$ws.Range("B11").Value = 1.135298427951739
$ws.Range("D11").Value = 0.01482462134583074
$ws.Range("E11").Value = 0.6935348852775078
$ws.Range("F11").Value = 0.4568281892317998
$ws.Range("G11").Value = 0.3168467701692634
$ws.Range("H11").Value = 0.386440763620115
$ws.Range("L11").Value = 0.2324773293739213
$ws.Range("M11").Value = 0.2628353816976556
$ws.Range("N11").Value = 1.854637232407271
$ws.Range("O11").Value = 1.333475205718884
# Row 12
$ws.Range("B12").Value = 1.147130303915418
$ws.Range("D12").Value = 0.01519536946602074
$ws.Range("E12").Value = 0.6950420686763792
$ws.Range("F12").Value = 0.4618205506234005
$ws.Range("G12").Value = 0.3214320167467122
$ws.Range("H12").Value = 0.3877692839377147
$ws.Range("L12").Value = 0.2375827454748816
$ws.Range("M12").Value = 0.2667650180801786
$ws.Range("N12").Value = 1.859774284043823
$ws.Range("O12").Value = 1.346043428948576
# Row 13
$ws.Range("B13").Value = 1.144580090165647
$ws.Range("D13").Value = 0.0151155365271407
$ws.Range("E13").Value = 0.6947165636007284
$ws.Range("F13").Value = 0.4607435182298616
$ws.Range("G13").Value = 0.320443179422071
$ws.Range("H13").Value = 0.3874817416817109
$ws.Range("L13").Value = 0.2364832843227589
$ws.Range("M13").Value = 0.265918420455705
$ws.Range("N13").Value = 1.85866355123855
$ws.Range("O13").Value = 1.343330925442757
# Row 14
$ws.Range("B14").Value = 1.136270937587454
$ws.Range("D14").Value = 0.01485512943220613
$ws.Range("E14").Value = 0.6936584731972104
$ws.Range("F14").Value = 0.4572380855640858
$ws.Range("G14").Value = 0.3172234053950973
$ws.Range("H14").Value = 0.3865494211250535
$ws.Range("L14").Value = 0.2328973913167971
$ws.Range("M14").Value = 0.2631585486404973
$ws.Range("N14").Value = 1.855057887025907
$ws.Range("O14").Value = 1.334506627401652
# Row 15
$ws.Range("B15").Value = 1.131187235008667
$ws.Range("D15").Value = 0.01469558102522228
$ws.Range("E15").Value = 0.6930130230077438
$ws.Range("F15").Value = 0.4550962886779928
$ws.Range("G15").Value = 0.3152550713778908
$ws.Range("H15").Value = 0.385982511099499
$ws.Range("L15").Value = 0.2307006949211967
$ws.Range("M15").Value = 0.2614688721378968
$ws.Range("N15").Value = 1.85286214617102
$ws.Range("O15").Value = 1.329118212131334
# Row 16
$ws.Range("B16").Value = 1.102145574330763
$ws.Range("D16").Value = 0.01378065155410724
$ws.Range("E16").Value = 0.6893553150419791
$ws.Range("F16").Value = 0.4429049965410314
$ws.Range("G16").Value = 0.3040345938116644
$ws.Range("H16").Value = 0.3827976118633671
$ws.Range("L16").Value = 0.2181089767294253
$ws.Range("M16").Value = 0.2517988549239405
$ws.Range("N16").Value = 1.840476666681241
$ws.Range("O16").Value = 1.298495566209937
# Row 17
$ws.Range("B17").Value = 1.084411355261835
$ws.Range("D17").Value = 0.01321883284421688
$ws.Range("E17").Value = 0.6871484624955784
$ws.Range("F17").Value = 0.4354996214862297
$ws.Range("G17").Value = 0.2972039798174393
$ws.Range("H17").Value = 0.380900762208114
$ws.Range("L17").Value = 0.2103816447098552
$ws.Range("M17").Value = 0.2458782789983616
$ws.Range("N17").Value = 1.833055542419004
$ws.Range("O17").Value = 1.279938029259512
# Row 18
$ws.Range("B18").Value = 1.074241481848048
$ws.Range("D18").Value = 0.01289549777328602
$ws.Range("E18").Value = 0.6858928888903506
$ws.Range("F18").Value = 0.4312672843776681
$ws.Range("G18").Value = 0.2932945564709968
$ws.Range("H18").Value = 0.3798306947329024
$ws.Range("L18").Value = 0.2059361905425448
$ws.Range("M18").Value = 0.2424772921189486
$ws.Range("N18").Value = 1.828852461276114
$ws.Range("O18").Value = 1.269348131244953
# Row 19
$ws.Range("B19").Value = 1.070803382134159
$ws.Range("D19").Value = 0.01278598972945844
$ws.Range("E19").Value = 0.6854701480632457
$ws.Range("F19").Value = 0.4298389312311315
$ws.Range("G19").Value = 0.2919742156919796
$ws.Range("H19").Value = 0.3794719873754104
$ws.Range("L19").Value = 0.2044308918012945
$ws.Range("M19").Value = 0.2413265352716536
$ws.Range("N19").Value = 1.82744062503383
$ws.Range("O19").Value = 1.265776979596666
# Row 20
$ws.Range("B20").Value = 1.086296055281196
$ws.Range("D20").Value = 0.01327865939923356
$ws.Range("E20").Value = 0.6873819658482034
$ws.Range("F20").Value = 0.4362851377859727
$ws.Range("G20").Value = 0.2979291053593869
$ws.Range("H20").Value = 0.3811005170244357
$ws.Range("L20").Value = 0.211204327565568
$ws.Range("M20").Value = 0.2465080838521772
$ws.Range("N20").Value = 1.833838777153289
$ws.Range("O20").Value = 1.281904825528585
# Row 21
$ws.Range("B21").Value = 1.13871031052679
$ws.Range("D21").Value = 0.01493162600111475
$ws.Range("E21").Value = 0.6939687062784543
$ws.Range("F21").Value = 0.4582665945197988
$ws.Range("G21").Value = 0.3181683242764706
$ws.Range("H21").Value = 0.3868223987104074
$ws.Range("L21").Value = 0.2339507042237585
$ws.Range("M21").Value = 0.2639690185000205
$ws.Range("N21").Value = 1.856114285464315
$ws.Range("O21").Value = 1.337095052488962
# Row 22
$ws.Range("B22").Value = 1.173230586653943
$ws.Range("D22").Value = 0.01601009285943178
$ws.Range("E22").Value = 0.698392986838634
$ws.Range("F22").Value = 0.4728736139411751
$ws.Range("G22").Value = 0.331568970955928
$ws.Range("H22").Value = 0.390748340777634
$ws.Range("L22").Value = 0.2488066065293992
$ws.Range("M22").Value = 0.2754178861953847
$ws.Range("N22").Value = 1.871247829840712
$ws.Range("O22").Value = 1.373913354146396
# Row 23
$ws.Range("B23").Value = 1.154782551442054
$ws.Range("D23").Value = 0.01543466989696896
$ws.Range("E23").Value = 0.696020877602102
$ws.Range("F23").Value = 0.4650555295986578
$ws.Range("G23").Value = 0.3244009171271216
$ws.Range("H23").Value = 0.3886359492100411
$ws.Range("L23").Value = 0.2408787678759268
$ws.Range("M23").Value = 0.2693040986776509
$ws.Range("N23").Value = 1.863118465576747
$ws.Range("O23").Value = 1.354194211331077
# Row 24
$ws.Range("B24").Value = 1.085443902077145
$ws.Range("D24").Value = 0.01325161286454346
$ws.Range("E24").Value = 0.6872763579110597
$ws.Range("F24").Value = 0.4359299276472512
$ws.Range("G24").Value = 0.2976012213261612
$ws.Range("H24").Value = 0.381010144147595
$ws.Range("L24").Value = 0.2108324016808467
$ws.Range("M24").Value = 0.2462233402390268
$ws.Range("N24").Value = 1.83348447910106
$ws.Range("O24").Value = 1.281015390627942
# Row 25
$ws.Range("B25").Value = 1.011815718223829
$ws.Range("D25").Value = 0.01088956952894193
$ws.Range("E25").Value = 0.678371389623841
$ws.Range("F25").Value = 0.4055467811216857
$ws.Range("G25").Value = 0.2694332153113237
$ws.Range("H25").Value = 0.3735855160792312
$ws.Range("L25").Value = 0.1783894096312082
$ws.Range("M25").Value = 0.2214951196794317
$ws.Range("N25").Value = 1.804021938070463
$ws.Range("O25").Value = 1.205287008960426
